# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" sheet right after "总计" (shifting 2022-Q2 / 2022-Q1 /
#    2021-Q4 / 2021-Q3 one slot to the right) and fill it with the new quarter's
#    fund-holdings detail table.
# 2) Insert a new row at the top of the "总计" summary sheet with the 2022-Q3
#    aggregate figures, pushing the older quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert the 2022-Q3 row
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 1.36

# Renumber the 0-based index column for the rows that shifted down
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# Part 2: new "2022-Q3" detail sheet
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q3.Name = "2022-Q3"

# Header row + first column share the "s=2" header/index style used by the
# other quarter sheets - copy it over from the summary sheet.
$summary.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A2").Copy()
$q3.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "002345", "华夏高端制造灵活配置混合A", "20.58", "90.97", "3.75", "0.7718", 9),
    @(1, "240004", "华宝动力组合混合A",         "14.13", "75.08", "2.76", "0.3900", 8),
    @(2, "016257", "华宝动力组合混合C",         "2.03",  "75.08", "2.76", "0.0560", 8),
    @(3, "003300", "华夏圆和灵活配置混合A",     "0.77",  "75.31", "5.69", "0.0438", 9),
    @(4, "011351", "金鹰年年邮益一年持有期混合A", "3.43",  "34.33", "1.15", "0.0394", 2),
    @(5, "015058", "华夏高端制造灵活配置混合C", "0.95",  "90.97", "3.75", "0.0356", 9),
    @(6, "015068", "华夏圆和灵活配置混合C",     "0.33",  "75.31", "5.69", "0.0188", 9),
    @(7, "011352", "金鹰年年邮益一年持有期混合C", "0.27",  "34.33", "1.15", "0.0031", 2)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $data[0]
    $q3.Cells.Item($r, 2).Value = "'" + $data[1]
    $q3.Cells.Item($r, 3).Value = $data[2]
    $q3.Cells.Item($r, 4).Value = "'" + $data[3]
    $q3.Cells.Item($r, 5).Value = "'" + $data[4]
    $q3.Cells.Item($r, 6).Value = "'" + $data[5]
    $q3.Cells.Item($r, 7).Value = "'" + $data[6]
    $q3.Cells.Item($r, 8).Value = $data[7]

    $q3.Range("B" + $r).Style = "Normal"
    $q3.Range("D" + $r + ":G" + $r).Style = "Normal"
}
